# Update coefficients with Carla
$wb = $excel.ActiveWorkbook

# --- Sheet "default" ---
$ws = $wb.Worksheets.Item("default")
$ws.Range("C3").Value = 0.4
$ws.Range("I3").Value = 0.00080000000000000004
$ws.Range("E5").Value = 0.54
$ws.Range("I5").Value = 0.00080000000000000004
[void]$ws.Range("C3").Select()

# --- Sheet "A" ---
$ws = $wb.Worksheets.Item("A")
$ws.Range("C3").Value = 0.4
$ws.Range("I3").Value = 0.00080000000000000004
$ws.Range("E5").Value = 0.54
$ws.Range("I5").Value = 0.00080000000000000004
[void]$ws.Range("C3").Select()

# --- Sheet "C10-C12" ---
$ws = $wb.Worksheets.Item("C10-C12")
$ws.Range("C3").Value = 0.4
$ws.Range("I3").Value = 0.00080000000000000004
$ws.Range("E5").Value = 0.54
$ws.Range("I5").Value = 0.00080000000000000004
[void]$ws.Range("C3").Select()

# --- Sheet "G-U_X_G4677" ---
$ws = $wb.Worksheets.Item("G-U_X_G4677")
$ws.Range("C3").Value = 0.4
$ws.Range("I3").Value = 0.00080000000000000004
$ws.Range("E5").Value = 0.54
$ws.Range("I5").Value = 0.00080000000000000004
[void]$ws.Range("E21").Select()

# --- Sheet "EP_HH" (selected/activated last to match the saved view state) ---
$ws = $wb.Worksheets.Item("EP_HH")
$ws.Range("C13").Value = 0.4
$ws.Range("E13").Value = 0.54
$ws.Range("I13").Value = 0.00051999999999999995
[void]$ws.Range("J20").Select()

Write-Host "done"
